$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.733.81"
$ws.Range("D3").Value = "2.075.05"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'232.90"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").Value = "'0.622"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "'58.37"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").Value = "'0.0782"
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("E11").Value = "  +3.06%  "
$ws.Range("D12").Value = "2.381.05"
$ws.Range("E12").Value = "  -1.41%  "
$ws.Range("D13").Value = "'14.73"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("D14").Value = "'20.82"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").Value = "'5.33"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("D17").Value = "2.075.07"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "37.690.14"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'6.17"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").Value = "'71.03"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("D21").Value = "0.0₃0832"
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("D22").Value = "'228.03"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'2.38"
$ws.Range("E24").Value = "  -1.52%  "
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("E27").Value = "  +4.32%  "
$ws.Range("D28").Value = "'9.01"
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("D29").Value = "'19.44"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("D31").Value = "'0.121"
$ws.Range("E31").Value = "  +2.46%  "
$ws.Range("D32").Value = "'4.67"
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("D33").Value = "'0.0630"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("D34").Value = "'4.65"
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'1.82"
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'3.39"
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "'5.31"
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'100.43"
$ws.Range("E40").Value = "  +4.04%  "
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").Value = "'0.0972"
$ws.Range("E41").Value = "  -3.57%  "
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("E43").Value = "  +1.00%  "
$ws.Range("D44").Value = "1.444.02"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("D45").Value = "'1.15"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("D46").Value = "'16.34"
$ws.Range("E46").Value = "  +6.37%  "
$ws.Range("D47").Value = "'4.19"
$ws.Range("E47").Value = "  +2.44%  "
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("D49").Value = "'7.40"
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").Value = "2.266.49"
$ws.Range("E51").Value = "  -1.44%  "
